# Generate Report for Handoff
# Adds two new tracked files (c824284b-... and ec5d911d-...) to the
# localization-status workbook: one row each on the "Overview" sheet and
# on the per-locale "zh-cn" / "de-de" sheets, mirroring the layout already
# used for the existing eb44f72a-... entry.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

$guid1 = "c824284b-6084-4a4f-9f03-0bda30486517"
$guid2 = "ec5d911d-ecf9-4ee6-bfc1-61e164c6eaaf"

$hash1 = "a1d16ea8beb9446fc096ef0c67854e15e0308290"
$hash2 = "7dbcd50e0d56501787c773a9c49a5ed9d0275aa3"

$status = "Ready for handoff"
$ext = ".md"
$reason = "Include"
$noHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A3").Value = "$guid1.md"
$wsOv.Range("B3").Value = $status
$wsOv.Range("C3").Value = $status
$wsOv.Range("D3").Value = "2016-03-22 10:40:10"
$wsOv.Range("D3").NumberFormat = $dateFmt
$wsOv.Hyperlinks.Add($wsOv.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f8247e2ae60c5b7c9a54f28be8899123793559c9/e2e/$guid1.md", "", "", "$guid1.md")
$wsOv.Range("A3").Style = "HyperLink"

$wsOv.Range("A4").Value = "$guid2.md"
$wsOv.Range("B4").Value = $status
$wsOv.Range("C4").Value = $status
$wsOv.Range("D4").Value = "2016-03-22 10:40:10"
$wsOv.Range("D4").NumberFormat = $dateFmt
$wsOv.Hyperlinks.Add($wsOv.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f8247e2ae60c5b7c9a54f28be8899123793559c9/e2e/$guid2.md", "", "", "$guid2.md")
$wsOv.Range("A4").Style = "HyperLink"

# ---------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest
# Handoff File | Latest Handoff Datetime | Latest Target File | Latest
# Handback File | Latest Handback DateTime | Reference Tokens |
# Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = "$guid1.md"
$wsZh.Range("B3").Value = $ext
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = "$guid1.$hash1.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-22 10:40:06"
$wsZh.Range("E3").NumberFormat = $dateFmt
$wsZh.Range("H3").Value = $noHandback
$wsZh.Range("H3").NumberFormat = $dateFmt
$wsZh.Range("J3").Value = $reason
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1b60366d052ae3fb97ee92283bc2e62def6cef79/e2e/$guid1.md", "", "", "$guid1.md")
$wsZh.Range("A3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/685a742235a8842be3b75fb92c45d2f1830e118a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$guid1.$hash1.zh-cn.xlf", "", "", "$guid1.$hash1.zh-cn.xlf")
$wsZh.Range("D3").Style = "HyperLink"

$wsZh.Range("A4").Value = "$guid2.md"
$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = "$guid2.$hash2.zh-cn.xlf"
$wsZh.Range("E4").Value = "2016-03-22 10:40:06"
$wsZh.Range("E4").NumberFormat = $dateFmt
$wsZh.Range("H4").Value = $noHandback
$wsZh.Range("H4").NumberFormat = $dateFmt
$wsZh.Range("J4").Value = $reason
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1b60366d052ae3fb97ee92283bc2e62def6cef79/e2e/$guid2.md", "", "", "$guid2.md")
$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/685a742235a8842be3b75fb92c45d2f1830e118a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$guid2.$hash2.zh-cn.xlf", "", "", "$guid2.$hash2.zh-cn.xlf")
$wsZh.Range("D4").Style = "HyperLink"

# ---------------------------------------------------------------------
# de-de sheet: same column layout as zh-cn
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = "$guid1.md"
$wsDe.Range("B3").Value = $ext
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = "$guid1.$hash1.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-22 10:40:10"
$wsDe.Range("E3").NumberFormat = $dateFmt
$wsDe.Range("H3").Value = $noHandback
$wsDe.Range("H3").NumberFormat = $dateFmt
$wsDe.Range("J3").Value = $reason
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b955305fc55ff3ef735c483dc2b619a3ad592118/e2e/$guid1.md", "", "", "$guid1.md")
$wsDe.Range("A3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/488fed9de7ca80adfa786b33031801f6b2fa1fa2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$guid1.$hash1.de-de.xlf", "", "", "$guid1.$hash1.de-de.xlf")
$wsDe.Range("D3").Style = "HyperLink"

$wsDe.Range("A4").Value = "$guid2.md"
$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = "$guid2.$hash2.de-de.xlf"
$wsDe.Range("E4").Value = "2016-03-22 10:40:10"
$wsDe.Range("E4").NumberFormat = $dateFmt
$wsDe.Range("H4").Value = $noHandback
$wsDe.Range("H4").NumberFormat = $dateFmt
$wsDe.Range("J4").Value = $reason
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b955305fc55ff3ef735c483dc2b619a3ad592118/e2e/$guid2.md", "", "", "$guid2.md")
$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/488fed9de7ca80adfa786b33031801f6b2fa1fa2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$guid2.$hash2.de-de.xlf", "", "", "$guid2.$hash2.de-de.xlf")
$wsDe.Range("D4").Style = "HyperLink"

Write-Output "Applied handoff rows for $guid1 and $guid2"
